$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 606, shifting existing rows 606..691 down to 607..692.
$ws.Rows.Item(606).Insert()

# Populate the newly inserted row 606 with its data.
$ws.Cells.Item(606, 1).Value = 6
$ws.Cells.Item(606, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(606, 3).Value = "Metropolitana"
$ws.Cells.Item(606, 4).Value = 44984
$ws.Cells.Item(606, 4).NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Cells.Item(606, 5).Value = 13
$ws.Cells.Item(606, 6).Value = 100112044
$ws.Cells.Item(606, 7).Value = "Perejil"
$ws.Cells.Item(606, 8).Value = "Sin especificar"
$ws.Cells.Item(606, 9).Value = "Primera"
$ws.Cells.Item(606, 10).Value = 240
$ws.Cells.Item(606, 11).Value = 12000
$ws.Cells.Item(606, 12).Value = 13000
$ws.Cells.Item(606, 13).Value = 12458
$ws.Cells.Item(606, 14).Value = "$/docena de atados"
$ws.Cells.Item(606, 15).Value = "Región Metropolitana"
$ws.Cells.Item(606, 16).Value = 4153
$ws.Cells.Item(606, 17).Value = 3
$ws.Cells.Item(606, 18).Value = "Hortaliza"
